# ============================================================================
# Edit: Settings sheet gets a new block of configuration rows (7-23), a
# hyperlink on B7, trailing blank rows removed, and the active sheet/
# selection moved. (Config.xlsx - "Se modifico el archivo de configuración")
# ============================================================================

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Settings")

# ----------------------------------------------------------------------
# 1) New configuration rows, columns A (key), B (value), C (description)
# ----------------------------------------------------------------------

$rows = @(
    @{ Row=7;  A="URL_SitioSupermercado";          B="https://www.laanonima.com.ar/";      C="URL base del portal comercial" },
    @{ Row=8;  A="Sucursal_BusquedaHabilitada";     B="TRUE";                                C="Indica si se debe seleccionar sucursal" },
    @{ Row=9;  A="Sucursal_CodigoPostal";           B=9420;                                  C="Código postal de la sucursal a seleccionar" },
    @{ Row=10; A="Sucursal_Denominacion";           B="RIO GRANDE 4";                        C="Denominación de la sucursal a seleccionar" },
    @{ Row=11; A="Sucursal_Domicilio";              B="AV.SAN MARTIN 1605";                  C="Domicilio de la sucursal a seleccionar" },
    @{ Row=12; A="Ruta_Input_Productos";            B="Data\Input\Productos_CBA.csv";        C="Archivo con el listado de productos CBA" },
    @{ Row=13; A="Ruta_Output_Resultados";          B="Data\Output\Resultados_CBA.csv";      C="Archivo CSV final generado por el robot" },
    @{ Row=14; A="Ruta_Output_RegistroTiempos";     B="Data\Output\RegistroTiempos.csv";     C="Registro de inicio/fin del proceso" },
    @{ Row=15; A="Ruta_Logs_Sesiones";              B="Data\Logs\LogSesiones\";              C="Carpeta para logs funcionales" },
    @{ Row=16; A="Ruta_Logs_Errores";               B="Data\Logs\LogErrores\";               C="Carpeta para logs técnicos" },
    @{ Row=17; A="Ruta_Screenshots";                B="Screenshots\";                        C="Capturas para errores y excepciones" },
    @{ Row=18; A="MaxRetriesBusiness";              B=0;                                     C="Negocio: no se reintenta" },
    @{ Row=19; A="MaxRetriesSystem";                B=2;                                     C="Reintentos para errores de aplicación" },
    @{ Row=20; A="TimeoutPagCarga";                 B=10000;                                 C="Timeout carga del sitio (ms)" },
    @{ Row=21; A="TimeoutBusqueda";                 B=5000;                                  C="Timeout búsqueda por producto (ms)" },
    @{ Row=22; A="TiempoEsperaResultados";          B=4000;                                  C="Tiempo de espera para lista de productos" },
    @{ Row=23; A="Correo_Notificaciones";           B="martinmirabete@gmail.com";            C="Dueño del proceso para alertas" }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.A
    $ws.Cells.Item($r.Row, 2).Value = $r.B
    $ws.Cells.Item($r.Row, 3).Value = $r.C
}

# ----------------------------------------------------------------------
# 2) Formatting: column A/C get a wrapped, vertically centered style;
#    column B (values) gets a wrapped, vertically centered style with a
#    smaller Arial Unicode MS font; B9 (postal code) is right aligned.
# ----------------------------------------------------------------------

# Base style for column A / C text cells (keys + descriptions)
$ws.Range("A8").Font.Name = "Calibri"
$ws.Range("A8").Font.Size = 11
$ws.Range("A8").WrapText = $true
$ws.Range("A8").VerticalAlignment = -4108
$ws.Range("A8").Copy()
$ws.Range("A7:A23").PasteSpecial(-4122)
$ws.Range("C7:C23").PasteSpecial(-4122)

# Base style for column B value cells (Arial Unicode MS 10pt)
$ws.Range("B8").Font.Name = "Arial Unicode MS"
$ws.Range("B8").Font.Size = 10
$ws.Range("B8").Font.Color = 0
$ws.Range("B8").WrapText = $true
$ws.Range("B8").VerticalAlignment = -4108
$ws.Range("B8").Copy()
$ws.Range("B8:B23").PasteSpecial(-4122)

# Right-aligned variant for the postal code value
$ws.Range("B9").HorizontalAlignment = -4152

# Empty, wrapped cell C6 (pre-existing wrap style reused)
$ws.Range("C6").WrapText = $true

$excel.CutCopyMode = $false

# ----------------------------------------------------------------------
# 3) Hyperlink on B7 (URL value) -- set value first so the hyperlink
#    creation doesn't clobber it with the display text.
# ----------------------------------------------------------------------

$ws.Range("B7").Value = "https://www.laanonima.com.ar/"
$ws.Range("B7").WrapText = $true
$ws.Range("B7").VerticalAlignment = -4108
$ws.Hyperlinks.Add($ws.Range("B7"), "https://www.laanonima.com.ar/", "", "", "https://www.supermercado.com.ar") | Out-Null
$ws.Range("B7").Value = "https://www.laanonima.com.ar/"

# ----------------------------------------------------------------------
# 4) Trim trailing blank rows 985-998 (dimension shrinks to A1:Z984)
# ----------------------------------------------------------------------

$ws.Rows("985:998").Delete()

# ----------------------------------------------------------------------
# 5) Sheet selection / active sheet: Settings becomes the active tab,
#    with the view scrolled down and a big selection block.
# ----------------------------------------------------------------------

$ws.Activate()
$ws.Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 3
$win.ScrollColumn = 1
$ws.Range("A24:XFD35").Select()

$wsAssets = $wb.Worksheets.Item("Assets")
$wsAssets.Range("A1").Select()

$ws.Activate()

$wb.Save()
